$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# Column A: 30.625 (raw) -> target 19.375; closest achievable via COM (MDW=7 pixel rounding) is 19.428571428571427
$ws.Columns.Item(1).ColumnWidth = 18.714285714285715
# Columns B:L: become a uniform narrow width (raw target 3.75); closest achievable is 3.714285714285714
$ws.Range("B1:L1").EntireColumn.ColumnWidth = 3.0

# --- New data rows (Miyachi VOICEROID / VOICEVOX skins) ---
# Values are assigned in the exact order needed so that the shared-string
# table indices line up with the authoritative edit (M70's text is added
# last, after all of the other new rows, matching the original edit order).
$ws.Range("A70").Value = "Miyachi"
$ws.Range("A71").Value = "Akane"
$ws.Range("M71").Value = "琴葉茜"
$ws.Range("A72").Value = "Akari"
$ws.Range("M72").Value = "紲星あかり"
$ws.Range("A73").Value = "Aoi"
$ws.Range("M73").Value = "琴葉葵"
$ws.Range("A74").Value = "Armachan"
$ws.Range("M74").Value = "アルマちゃん"
$ws.Range("A75").Value = "Diachan"
$ws.Range("M75").Value = "ディアちゃん"
$ws.Range("A76").Value = "flower"
$ws.Range("M76").Value = "flower"
$ws.Range("A77").Value = "Hau"
$ws.Range("M77").Value = "雨晴はう"
$ws.Range("A78").Value = "Hime"
$ws.Range("M78").Value = "鳴花ヒメ"
$ws.Range("A79").Value = "KANATA"
$ws.Range("M79").Value = "カナタ"
$ws.Range("A80").Value = "Kotoe"
$ws.Range("M80").Value = "タンゲコトエ"
$ws.Range("A81").Value = "Kou"
$ws.Range("M81").Value = "水奈瀬コウ"
$ws.Range("A82").Value = "MANA"
$ws.Range("M82").Value = "MANA"
$ws.Range("A83").Value = "Metan"
$ws.Range("M83").Value = "四国めたん"
$ws.Range("A84").Value = "Mikoto"
$ws.Range("M84").Value = "鳴花ミコト"
$ws.Range("A85").Value = "NAKO"
$ws.Range("M85").Value = "ナコ"
$ws.Range("A86").Value = "Rei"
$ws.Range("M86").Value = "足立レイ"
$ws.Range("A87").Value = "REKO"
$ws.Range("M87").Value = "レコ"
$ws.Range("A88").Value = "Ritsu"
$ws.Range("M88").Value = "波音リツ"
$ws.Range("A89").Value = "Rowen"
$ws.Range("M89").Value = "式狼縁"
$ws.Range("A90").Value = "Sora"
$ws.Range("M90").Value = "桜乃そら"
$ws.Range("A91").Value = "Taigen"
$ws.Range("M91").Value = "式大元"
$ws.Range("A92").Value = "Tsukuyomichan"
$ws.Range("M92").Value = "つくよみちゃん"
$ws.Range("A93").Value = "Tsumugi"
$ws.Range("M93").Value = "春日部つむぎ"
$ws.Range("A94").Value = "Yukari"
$ws.Range("M94").Value = "結月ゆかり"
$ws.Range("A95").Value = "Yuzuru"
$ws.Range("M95").Value = "伊織弓鶴"
$ws.Range("M96").Value = "ずんだもん"
$ws.Range("A96").Value = "Zundamon"
$ws.Range("M70").Value = "みやち作"

# --- Final selection state ---
$ws.Range("M96").Select()
